# Add the "2023" column (T) to the Hepatitis B incidence table on the
# single worksheet: copy the formatting of column S (years 2007-2022)
# into the new column T, then populate the 2023 header + data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the per-row cell styles from S3:S33 into T3:T33 so each new cell
# picks up exactly the same number format / font / border as its row.
$ws.Range("S3:S33").Copy()
$ws.Range("T3:T33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header: year 2023
$ws.Range("T3").Value = 2023

# Data values for 2023, one per indicator row
$ws.Range("T4").Value = 2.3381104968484805
$ws.Range("T5").Value = 2.0344672190198714
$ws.Range("T6").Value = 2.6483752218014245
$ws.Range("T7").Value = 3.9852372948902328
$ws.Range("T8").Value = 4.5532396299967433
$ws.Range("T9").Value = 3.4291318466903733
$ws.Range("T10").Value = 1.2089851778417198
$ws.Range("T11").Value = 1.521116134174612
$ws.Range("T12").Value = 0.9008846687447073
$ws.Range("T13").Value = 3.694303753043183
$ws.Range("T14").Value = 4.0607488020791038
$ws.Range("T15").Value = 3.327319511401615
$ws.Range("T16").Value = 0.32236434908190637
$ws.Range("T17").Value = 0
$ws.Range("T18").Value = 0.63756806039044667
$ws.Range("T19").Value = 2.1691385808410835
$ws.Range("T20").Value = 1.5024572004578396
$ws.Range("T21").Value = 2.8259763748375066
$ws.Range("T22").Value = 6.1744985943935555
$ws.Range("T23").Value = 4.3993752887090034
$ws.Range("T24").Value = 7.9169155696940479
$ws.Range("T25").Value = 2.8763040791558883
$ws.Range("T26").Value = 1.4751329463567904
$ws.Range("T27").Value = 4.2954684675262591
$ws.Range("T28").Value = 1.8177568880002077
$ws.Range("T29").Value = 1.581380197008345
$ws.Range("T30").Value = 2.103608453446189
$ws.Range("T31").Value = 1.3736037318066185
$ws.Range("T32").Value = 2.249820014398848
$ws.Range("T33").Value = 0.53701655085009725

# Row 4 picked up an explicit custom height in the authored edit.
$ws.Rows.Item(4).RowHeight = 16.5
